$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Geography column (D) values per refined taxonomy
$ws.Range("D2").Value = "Global/Multi-region"
$ws.Range("D4").Value = "Asia - Emerging"
$ws.Range("D6").Value = "Latin America"
$ws.Range("D7").Value = "Sub-Saharan Africa"
$ws.Range("D8").Value = "Global/Multi-region"
$ws.Range("D10").Value = "Asia - Emerging"
$ws.Range("D12").Value = "Asia - Emerging"
$ws.Range("D14").Value = "Asia - Developed"
$ws.Range("D15").Value = "Asia - Emerging"
$ws.Range("D16").Value = "Europe (excluding UK)"
$ws.Range("D17").Value = "UK"
$ws.Range("D18").Value = "Europe (excluding UK)"
$ws.Range("D19").Value = "Asia - Emerging"
$ws.Range("D20").Value = "Australia/Oceania"
$ws.Range("D21").Value = "Europe (excluding UK)"
$ws.Range("D23").Value = "Asia - Emerging"
$ws.Range("D25").Value = "Asia - Emerging"
$ws.Range("D26").Value = "Middle East & North Africa"
$ws.Range("D27").Value = "Europe (excluding UK)"
$ws.Range("D28").Value = "Sub-Saharan Africa"
$ws.Range("D29").Value = "Asia - Emerging"
$ws.Range("D31").Value = "Asia - Developed"
$ws.Range("D32").Value = "Europe (excluding UK)"
$ws.Range("D34").Value = "UK"
$ws.Range("D35").Value = "UK"
$ws.Range("D36").Value = "UK"
$ws.Range("D37").Value = "UK"
$ws.Range("D38").Value = "UK"
$ws.Range("D39").Value = "Europe (excluding UK)"
$ws.Range("D40").Value = "Europe (excluding UK)"
$ws.Range("D41").Value = "Europe (excluding UK)"
$ws.Range("D42").Value = "UK"
$ws.Range("D43").Value = "Europe (excluding UK)"

# Narrow column D to fit the new shorter labels.
# The engine stores XML column width as (ColumnWidth + 0.8333...), matching
# Excel's character-width -> internal-width padding quirk, so back the
# offset out to land exactly on the target stored width of 28.
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668
